$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update each attendance row: date label (DD-MM-YYYY text) and D..H counts.
# NumberFormat "@" forces the date-like text to be stored as a literal string
# instead of being auto-parsed into a date serial value; the style is then
# reset back to Normal so no stray formatting is left on the cell.
$a = $ws.Range("A3")
$a.NumberFormat = "@"
$a.Value = "28-07-2022"
$a.Style = "Normal"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$a = $ws.Range("A4")
$a.NumberFormat = "@"
$a.Value = "01-08-2022"
$a.Style = "Normal"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

$a = $ws.Range("A5")
$a.NumberFormat = "@"
$a.Value = "04-08-2022"
$a.Style = "Normal"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

$a = $ws.Range("A6")
$a.NumberFormat = "@"
$a.Value = "08-08-2022"
$a.Style = "Normal"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

$a = $ws.Range("A7")
$a.NumberFormat = "@"
$a.Value = "11-08-2022"
$a.Style = "Normal"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1

$a = $ws.Range("A8")
$a.NumberFormat = "@"
$a.Value = "15-08-2022"
$a.Style = "Normal"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1

$a = $ws.Range("A9")
$a.NumberFormat = "@"
$a.Value = "18-08-2022"
$a.Style = "Normal"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 1

$a = $ws.Range("A10")
$a.NumberFormat = "@"
$a.Value = "22-08-2022"
$a.Style = "Normal"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1

$a = $ws.Range("A11")
$a.NumberFormat = "@"
$a.Value = "25-08-2022"
$a.Style = "Normal"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0

$a = $ws.Range("A12")
$a.NumberFormat = "@"
$a.Value = "29-08-2022"
$a.Style = "Normal"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0

$a = $ws.Range("A13")
$a.NumberFormat = "@"
$a.Value = "01-09-2022"
$a.Style = "Normal"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0

$a = $ws.Range("A14")
$a.NumberFormat = "@"
$a.Value = "05-09-2022"
$a.Style = "Normal"
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 1

$a = $ws.Range("A15")
$a.NumberFormat = "@"
$a.Value = "08-09-2022"
$a.Style = "Normal"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 1

$a = $ws.Range("A16")
$a.NumberFormat = "@"
$a.Value = "12-09-2022"
$a.Style = "Normal"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0

$a = $ws.Range("A17")
$a.NumberFormat = "@"
$a.Value = "15-09-2022"
$a.Style = "Normal"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1

$a = $ws.Range("A18")
$a.NumberFormat = "@"
$a.Value = "19-09-2022"
$a.Style = "Normal"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 1

$a = $ws.Range("A19")
$a.NumberFormat = "@"
$a.Value = "22-09-2022"
$a.Style = "Normal"
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 1

$a = $ws.Range("A20")
$a.NumberFormat = "@"
$a.Value = "26-09-2022"
$a.Style = "Normal"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0

$a = $ws.Range("A21")
$a.NumberFormat = "@"
$a.Value = "29-09-2022"
$a.Style = "Normal"
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 1
